$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text representation
# (values like "1.002" or "28.834.97" must not be auto-converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.834.97"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.880.62"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "325.21"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "0.4593"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "0.3884"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "0.07859"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "0.9849"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").Value = "21.75"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.870.39"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "7.004"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "5.666"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "0.06940"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "88.24"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "0.000009956"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "16.97"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "28.862.28"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "5.269"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").Value = "10.94"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").Value = "2.159.70"
$ws.Range("E24").Value = "  +2.19%  "
$ws.Range("D25").Value = "2.089"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "155.47"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "19.29"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "5.954"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("D29").Value = "1.929"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").Value = "117.49"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").Value = "0.09339"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").Value = "0.9016"
$ws.Range("E32").Value = "  -2.46%  "
$ws.Range("D33").Value = "5.269"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").Value = "1.327"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").Value = "3.266"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "1.191"
$ws.Range("E36").Value = "  +2.57%  "
$ws.Range("D37").Value = "0.05760"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "0.02072"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").Value = "1.002"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "7.658"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").Value = "0.5658"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").Value = "0.1765"
$ws.Range("E42").Value = "  -1.18%  "
$ws.Range("D43").Value = "9.687"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "2.255"
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").Value = "11.93"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("D46").Value = "0.5349"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").Value = "0.07039"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").Value = "1.847"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").Value = "2.537"
$ws.Range("E49").Value = "  +2.94%  "
$ws.Range("D50").Value = "112.75"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  -4.91%  "
